$p = $ppt.ActivePresentation
$s2 = $p.Slides.Item(2)
$grp52 = $s2.Shapes.Item(1)
$ungrouped = $grp52.Ungroup()

$g43 = $null
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $it = $s2.Shapes.Item($i)
    if ($it.Name -eq "Group 43") { $g43 = $it }
}
$ungrouped43 = $g43.Ungroup()

for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $it = $s2.Shapes.Item($i)
    if ($it.Name -eq "Picture 11") {
        $it.Delete()
        break
    }
}
Write-Output ("Top shapes after delete: " + $s2.Shapes.Count)

$names = @()
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $names += $s2.Shapes.Item($i).Name
}
Write-Output ("names: " + ($names -join ","))
$range = $s2.Shapes.Range($names)
Write-Output ("Range count: " + $range.Count)
$regrouped = $range.Group()
Write-Output ("Regrouped name: [" + $regrouped.Name + "]")
Write-Output ("Top shapes after regroup: " + $s2.Shapes.Count)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $it = $s2.Shapes.Item($i)
    Write-Output ("  shape $i : " + $it.Name + " id=" + $it.Id)
}
